$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.569.40'
Set-TextValue $ws.Range('E2') '  -0.60%  '
Set-TextValue $ws.Range('D3') '1.836.24'
Set-TextValue $ws.Range('E3') '  -0.74%  '
Set-TextValue $ws.Range('E4') '  -0.02%  '
Set-TextValue $ws.Range('D5') '314.33'
Set-TextValue $ws.Range('E5') '  +0.08%  '
Set-TextValue $ws.Range('E6') '  -0.03%  '
Set-TextValue $ws.Range('D7') '0.4291'
Set-TextValue $ws.Range('E7') '  -0.93%  '
Set-TextValue $ws.Range('D8') '0.3663'
Set-TextValue $ws.Range('E8') '  +0.16%  '
Set-TextValue $ws.Range('B9') 'Dogecoin'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D9') '0.07279'
Set-TextValue $ws.Range('E9') '  -0.73%  '
Set-TextValue $ws.Range('B10') 'Polygon'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D10') '0.8720'
Set-TextValue $ws.Range('E10') '  -0.77%  '
Set-TextValue $ws.Range('B11') 'Solana'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range('D11') '20.74'
Set-TextValue $ws.Range('E11') '  -0.07%  '
Set-TextValue $ws.Range('B12') 'WrappedEther'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D12') '1.773.31'
Set-TextValue $ws.Range('E12') '  -1.76%  '
Set-TextValue $ws.Range('B13') 'Polkadot'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D13') '5.436'
Set-TextValue $ws.Range('E13') '  +1.66%  '
Set-TextValue $ws.Range('B14') 'Chainlink'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D14') '6.545'
Set-TextValue $ws.Range('E14') '  +0.22%  '
Set-TextValue $ws.Range('B15') 'TRON'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D15') '0.06935'
Set-TextValue $ws.Range('E15') '  +0.04%  '
Set-TextValue $ws.Range('B16') 'BinanceUSD'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D16') '1.004'
Set-TextValue $ws.Range('B17') 'Litecoin'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D17') '80.40'
Set-TextValue $ws.Range('E17') '  +0.13%  '
Set-TextValue $ws.Range('B18') 'ShibaInu'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D18') '0.000008933'
Set-TextValue $ws.Range('E18') '  -1.18%  '
Set-TextValue $ws.Range('B19') 'Dai'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D19') '1.001'
Set-TextValue $ws.Range('E19') '  +0.05%  '
Set-TextValue $ws.Range('B20') 'Avalanche'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D20') '15.43'
Set-TextValue $ws.Range('E20') '  +0.15%  '
Set-TextValue $ws.Range('B21') 'WrappedBTC'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D21') '27.336.35'
Set-TextValue $ws.Range('E21') '  -1.02%  '
Set-TextValue $ws.Range('B22') 'Uniswap'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D22') '5.173'
Set-TextValue $ws.Range('E22') '  +3.79%  '
Set-TextValue $ws.Range('B23') 'Cosmos'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D23') '10.88'
Set-TextValue $ws.Range('E23') '  +4.92%  '
Set-TextValue $ws.Range('B24') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D24') '1.992.42'
Set-TextValue $ws.Range('E24') '  -2.91%  '
Set-TextValue $ws.Range('B25') 'Toncoin'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D25') '1.979'
Set-TextValue $ws.Range('E25') '  -0.39%  '
Set-TextValue $ws.Range('B26') 'Monero'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D26') '154.73'
Set-TextValue $ws.Range('E26') '  -0.56%  '
Set-TextValue $ws.Range('B27') 'EthereumClassic'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '18.93'
Set-TextValue $ws.Range('E27') '  +1.71%  '
Set-TextValue $ws.Range('B28') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D28') '5.185'
Set-TextValue $ws.Range('E28') '  -1.66%  '
Set-TextValue $ws.Range('B29') 'BitcoinCash'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D29') '114.81'
Set-TextValue $ws.Range('E29') '  -4.97%  '
Set-TextValue $ws.Range('B30') 'LidoDAOToken'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D30') '1.838'
Set-TextValue $ws.Range('E30') '  -1.57%  '
Set-TextValue $ws.Range('B31') 'Stellar'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D31') '0.08908'
Set-TextValue $ws.Range('E31') '  -0.22%  '
Set-TextValue $ws.Range('B32') 'ImmutableX'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D32') '0.7599'
Set-TextValue $ws.Range('E32') '  +0.37%  '
Set-TextValue $ws.Range('B33') 'Filecoin'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D33') '4.551'
Set-TextValue $ws.Range('E33') '  +0.18%  '
Set-TextValue $ws.Range('B34') 'HuobiToken'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D34') '2.964'
Set-TextValue $ws.Range('E34') '  +0.35%  '
Set-TextValue $ws.Range('B35') 'ARBITRUM'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D35') '1.144'
Set-TextValue $ws.Range('E35') '  +2.08%  '
Set-TextValue $ws.Range('B36') 'Frax'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D36') '1.001'
Set-TextValue $ws.Range('E36') '  -0.01%  '
Set-TextValue $ws.Range('B37') 'TrustWalletToken'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D37') '1.096'
Set-TextValue $ws.Range('E37') '  -1.41%  '
Set-TextValue $ws.Range('B38') 'Hedera'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D38') '0.05325'
Set-TextValue $ws.Range('E38') '  -1.66%  '
Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.01947'
Set-TextValue $ws.Range('E39') '  +0.55%  '
Set-TextValue $ws.Range('B40') 'MXToken'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D40') '2.810'
Set-TextValue $ws.Range('E40') '  -0.79%  '
Set-TextValue $ws.Range('B41') 'Algorand'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D41') '0.1673'
Set-TextValue $ws.Range('E41') '  +0.76%  '
Set-TextValue $ws.Range('B42') 'TheSandbox'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D42') '0.5102'
Set-TextValue $ws.Range('E42') '  +0.12%  '
Set-TextValue $ws.Range('B43') 'FraxShare'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '6.608'
Set-TextValue $ws.Range('E43') '  -0.87%  '
Set-TextValue $ws.Range('B44') 'Aptos'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D44') '8.447'
Set-TextValue $ws.Range('E44') '  +1.44%  '
Set-TextValue $ws.Range('B45') 'EnergySwap'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D45') '10.53'
Set-TextValue $ws.Range('E45') '  +1.90%  '
Set-TextValue $ws.Range('B46') 'Quant'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D46') '106.31'
Set-TextValue $ws.Range('E46') '  +1.88%  '
Set-TextValue $ws.Range('B47') 'Cronos'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D47') '0.06506'
Set-TextValue $ws.Range('E47') '  -0.52%  '
Set-TextValue $ws.Range('B48') 'Decentraland'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D48') '0.4685'
Set-TextValue $ws.Range('E48') '  +0.43%  '
Set-TextValue $ws.Range('B49') 'PaxDollar'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D49') '1.001'
Set-TextValue $ws.Range('E49') '  +0.01%  '
Set-TextValue $ws.Range('B50') 'NEARProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D50') '1.622'
Set-TextValue $ws.Range('E50') '  -0.12%  '
Set-TextValue $ws.Range('B51') 'RenderToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D51') '1.761'
Set-TextValue $ws.Range('E51') '  +2.82%  '
